# EmailNoMatch.xlsx - update execution timestamps for the VLink "Card Not
# Accepted" run (rows 2-7, column B) to the newest recorded test pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(
    "Wed Feb 08 16:40:16 EST 2023",
    "Wed Feb 08 16:40:27 EST 2023",
    "Wed Feb 08 16:40:38 EST 2023",
    "Wed Feb 08 16:40:48 EST 2023",
    "Wed Feb 08 16:40:59 EST 2023",
    "Wed Feb 08 16:41:11 EST 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
